$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (last changed) date serial.
# Every data row (2..lastRow) needs its value bumped from 46075 to 46076,
# i.e. incremented by one day.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
